$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Stable "donor" cells used purely as style/type templates. None of these
# cells (row 14) are themselves modified by this script, so they stay valid
# templates throughout the whole run.
#   C14 -> style 14, shared text "0"
#   E14 -> style 14, shared text "***.*"
#   I14 -> style 15, plain integer count cell
#   N14 -> style 16, plain decimal/percent cell
# ---------------------------------------------------------------------------
$donorText0    = $ws.Range("C14")
$donorTextStar = $ws.Range("E14")
$donorNum15    = $ws.Range("I14")
$donorNum16    = $ws.Range("N14")

# ---------------------------------------------------------------------------
# Header text updates (Volume/Number + report week dates)
# ---------------------------------------------------------------------------
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "29"

$weekCell = $ws.Range("C9")
$weekCell.Characters(27, 9).Text = "7/17/2023"
$weekCell.Characters(47, 9).Text = "7/23/2023"

# ---------------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------------
$donorText0.Copy($ws.Range("D16"))
$donorTextStar.Copy($ws.Range("E16"))
$ws.Range("F16").Value = 2
$ws.Range("H16").Value = 0
$ws.Range("L16").Value = 120

# ---------------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 2
$donorText0.Copy($ws.Range("D17"))
$donorTextStar.Copy($ws.Range("E17"))
$ws.Range("F17").Value = 11
$ws.Range("H17").Value = 57.142857142857
$ws.Range("I17").Value = 56
$ws.Range("K17").Value = 143.478260869565
$ws.Range("L17").Value = 166.666666666667
$ws.Range("M17").Value = 143.478260869565
$ws.Range("N17").Value = -6.666666666666

# ---------------------------------------------------------------------------
# Row 18
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 2
$donorNum15.Copy($ws.Range("D18"))
$ws.Range("D18").Value = 1
$donorNum16.Copy($ws.Range("E18"))
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 7
$ws.Range("H18").Value = 75
$ws.Range("I18").Value = 36
$ws.Range("J18").Value = 18
$ws.Range("K18").Value = 100
$ws.Range("L18").Value = 125
$ws.Range("M18").Value = -41.935483870967
$ws.Range("N18").Value = -78.313253012048

# ---------------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -12.5
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = -7.692307692307
$ws.Range("I19").Value = 165
$ws.Range("J19").Value = 154
$ws.Range("K19").Value = 7.142857142857
$ws.Range("L19").Value = 83.333333333333
$ws.Range("M19").Value = 126.027397260274
$ws.Range("N19").Value = 39.830508474576

# ---------------------------------------------------------------------------
# Row 20
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 5
$donorText0.Copy($ws.Range("D20"))
$donorTextStar.Copy($ws.Range("E20"))
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 233.333333333333
$ws.Range("I20").Value = 50
$ws.Range("K20").Value = -1.960784313725
$ws.Range("L20").Value = 117.391304347826
$ws.Range("M20").Value = 127.272727272727
$ws.Range("N20").Value = -86.945169712793

# ---------------------------------------------------------------------------
# Row 21 (TOTAL row - styles unchanged, values only)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 9
$ws.Range("E21").Value = 77.777777777777
$ws.Range("F21").Value = 54
$ws.Range("G21").Value = 42
$ws.Range("H21").Value = 28.571428571428
$ws.Range("I21").Value = 320
$ws.Range("J21").Value = 254
$ws.Range("K21").Value = 25.984251968503
$ws.Range("L21").Value = 103.821656050955
$ws.Range("M21").Value = 64.102564102564
$ws.Range("N21").Value = -57.839262187088

# ---------------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = -50
$ws.Range("F24").Value = 38
$ws.Range("G24").Value = 30
$ws.Range("H24").Value = 26.666666666666
$ws.Range("I24").Value = 272
$ws.Range("J24").Value = 227
$ws.Range("K24").Value = 19.823788546255
$ws.Range("L24").Value = 97.101449275362
$ws.Range("M24").Value = -5.226480836236

# ---------------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 6
$donorNum15.Copy($ws.Range("D25"))
$ws.Range("D25").Value = 3
$donorNum16.Copy($ws.Range("E25"))
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 18
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 38.461538461538
$ws.Range("I25").Value = 111
$ws.Range("J25").Value = 96
$ws.Range("K25").Value = 15.625
$ws.Range("L25").Value = 50
$ws.Range("M25").Value = -5.128205128205

# ---------------------------------------------------------------------------
# Row 27
# ---------------------------------------------------------------------------
$donorNum15.Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$donorNum15.Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$donorNum16.Copy($ws.Range("E27"))
$ws.Range("E27").Value = 0
$donorNum15.Copy($ws.Range("F27"))
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 11
$ws.Range("J27").Value = 11
$ws.Range("L27").Value = 37.5

# ---------------------------------------------------------------------------
# Row 28
# ---------------------------------------------------------------------------
$donorText0.Copy($ws.Range("C28"))
$donorText0.Copy($ws.Range("D28"))
$donorTextStar.Copy($ws.Range("E28"))

# ---------------------------------------------------------------------------
# Row 29
# ---------------------------------------------------------------------------
$donorText0.Copy($ws.Range("C29"))
$donorText0.Copy($ws.Range("D29"))
$donorTextStar.Copy($ws.Range("E29"))
